$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ClearCanvas toolkit licence was actually GPLv3, not GNU.
$ws.Range("C2").Value = "GPLv3"

# Add a new licence entry for the "Now UI Dashboard" CSS template that is
# used under WebService/assets, copying the formatting of the preceding
# (MathNet.Numerics) row so the new row matches the rest of the table.
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)

# Give the link its own hyperlink (this becomes rId1 on the sheet).
$ws.Hyperlinks.Add($ws.Range("D11"), "https://www.creative-tim.com/product/now-ui-dashboard")

$ws.Range("C11").Value = "MIT"
$ws.Range("A11").Value = "Now UI Dashboard"
$ws.Range("B11").Value = "WebService/assets"
$ws.Range("E11").Value = "CSS template for web service"

# Turn the existing ClearCanvas link into a real hyperlink too (rId2).
$ws.Hyperlinks.Add($ws.Range("D2"), "https://clearcanvas.github.io/")

# Restore the selection to C2, as left by the author after editing the
# licence cell.
$ws.Range("C2").Select()
